# Daily attendance processing - 2026-01-03 15:33:04
#
# The "Recorded By" column (G) holds a comma-separated list of the
# users/systems that touched each attendance row. This pass normalizes
# that list by moving the last contributor to the front (a right
# rotation by one position) for every multi-contributor cell - except
# the two known system-only combinations that must stay as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combinations that are intentionally left untouched by this pass.
$exceptions = @("admin@admin.com, System", "System, backup@backdoor.com")

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($exceptions -contains $val) { continue }

    $parts = $val.Split(",")
    if ($parts.Count -le 1) { continue }

    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $n = $trimmed.Count
    $lastEntry = $trimmed[$n - 1]
    $rest = $trimmed[0..($n - 2)]
    $rotated = @($lastEntry) + $rest

    $cell.Value = [string]::Join(", ", $rotated)
}
